$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Activate()

# Insert a new blank column before column N (shifts Late/heading/Outstanding
# columns one to the right), matching the formatting of the column to its
# left (M) for the new column's width.
$mWidth = $ws.Columns("M").ColumnWidth
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $mWidth

# Move the active selection to K18 on the Repayment schedule sheet (which
# also becomes the active/selected tab).
$ws.Range("K18").Select() | Out-Null
